# Auto-generated edit script applying cryptos list price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.847.67"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "2.439.03"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'560.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").Value = "'162.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("D9").Value = "'0.169"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.55%  "

$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "'4.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.74%  "

$ws.Range("E13").Value = "  +5.92%  "

$ws.Range("D14").Value = "68.734.33"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "2.887.44"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "2.443.51"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").Value = "'10.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").Value = "'339.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "'6.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").Value = "'3.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "

$ws.Range("E22").Value = "  +2.85%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'67.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("D26").Value = "2.568.31"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "'8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").Value = "0.0₃0822"
$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("D30").Value = "'7.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'428.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "

$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").Value = "'158.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "'19.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").Value = "'17.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "

$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("E41").Value = "  +4.12%  "

$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("E44").Value = "  +0.79%  "

$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").Value = "'130.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").Value = "'0.0719"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").Value = "'0.556"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "'0.0925"
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = "  +1.27%  "
